$p = $ppt.ActivePresentation
$nm = $p.NotesMaster
$theme = $nm.Theme
$tcs = $theme.ThemeColorScheme
$c1 = $tcs.Colors(1)
$c1.RGB = 255
